$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'247.21"
$ws.Range("D2").Style = "Normal"
$ws.Range("D3").Value = "'26.45"
$ws.Range("D3").Style = "Normal"
$ws.Range("D4").Value = "'5.083"
$ws.Range("D4").Style = "Normal"
$ws.Range("D5").Value = "'0.05618"
$ws.Range("D5").Style = "Normal"
$ws.Range("D6").Value = "'6.519"
$ws.Range("D6").Style = "Normal"
$ws.Range("D7").Value = "'0.8148"
$ws.Range("D7").Style = "Normal"
$ws.Range("D8").Value = "'0.8483"
$ws.Range("D8").Style = "Normal"
$ws.Range("B9").Value = "LiechtensteinCryptoassetsExchange"
$ws.Range("C9").Value = "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
$ws.Range("D9").Value = "'0.03193"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "8LiechtensteinCryptoassetsExchangeLCX"
$ws.Range("B10").Value = "BitrueCoin"
$ws.Range("C10").Value = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
$ws.Range("D10").Value = "'0.02840"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "9BitrueCoinBTR"
$ws.Range("B11").Value = "BitMartToken"
$ws.Range("C11").Value = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
$ws.Range("D11").Value = "'0.09401"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "10BitMartTokenBMX"
$ws.Range("B12").Value = "BitForexToken"
$ws.Range("C12").Value = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
$ws.Range("D12").Value = "'0.001523"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "11BitForexTokenBF"
$ws.Range("D13").Value = "'0.006181"
$ws.Range("D13").Style = "Normal"
$ws.Range("D14").Value = "'3.587"
$ws.Range("D14").Style = "Normal"
$ws.Range("D15").Value = "'3.059"
$ws.Range("D15").Style = "Normal"
$ws.Range("D18").Value = "'0.1347"
$ws.Range("D18").Style = "Normal"
$ws.Range("D19").Value = "'0.06992"
$ws.Range("D19").Style = "Normal"
$ws.Range("B20").Value = "ProBitToken"
$ws.Range("C20").Value = "https://coinranking.com/coin/lQP4d6T2+probittoken-prob"
$ws.Range("D20").Value = "'0.1320"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "19ProBitTokenPROB"
$ws.Range("B21").Value = "MCDex"
$ws.Range("C21").Value = "https://coinranking.com/coin/3nMM61qeg+mcdex-mcb"
$ws.Range("D21").Value = "'3.742"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "20MCDexMCB"
$ws.Range("B22").Value = "CoinExToken"
$ws.Range("C22").Value = "https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet"
$ws.Range("D22").Value = "'0.04673"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "21CoinExTokenCET"
$ws.Range("B23").Value = "ZBToken"
$ws.Range("C23").Value = "https://coinranking.com/coin/CxmvOsCyENPso+zbtoken-zb"
$ws.Range("D23").Value = "'0.1350"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "22ZBTokenZB"
$ws.Range("B24").Value = "One"
$ws.Range("C24").Value = "https://coinranking.com/coin/6Lga5NiXX3rT+one-one"
$ws.Range("D24").Value = "'0.0005970"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "23OneONE"
$ws.Range("D25").Value = "'0.001251"
$ws.Range("D25").Style = "Normal"
$ws.Range("D26").Value = "'0.004615"
$ws.Range("D26").Style = "Normal"
$ws.Range("D27").Value = "'0.00009601"
$ws.Range("D27").Style = "Normal"
$ws.Range("D40").Value = "'0.03674"
$ws.Range("D40").Style = "Normal"
$ws.Range("B41").Value = "BKEXToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk"
$ws.Range("D41").Value = "'0.1357"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "40BKEXTokenBKKBestin24h"
$ws.Range("B42").Value = "KickToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick"
$ws.Range("D42").Value = "'0.006121"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "41KickTokenKICK"
$ws.Range("D43").Value = "'0.002660"
$ws.Range("D43").Style = "Normal"
$ws.Range("D44").Value = "'0.008588"
$ws.Range("D44").Style = "Normal"
$ws.Range("D45").Value = "'0.00005293"
$ws.Range("D45").Style = "Normal"
$ws.Range("D48").Value = "'0.002277"
$ws.Range("D48").Style = "Normal"
